# Update the Daily APR data sheet:
#  - correct the redemption_rate value recorded for row id=3 (cell C4)
#  - append a new snapshot row (id=4) with timestamp + redemption_rate

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the previously-recorded redemption rate for the 3rd snapshot.
$ws.Range("C4").Value = 1.6694299880323946

# Append the new 4th snapshot row.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "2025-08-31T13:39"
$ws.Range("C5").Value = 1.6703984340747713
